$d = $word.ActiveDocument

# 1) The empty paragraph after "Add new obstacles based off turns..." becomes
#    "Make Bat predict rat position" (with gramStart/gramEnd proofErr markers
#    wrapping the final word, mirroring the surrounding document's style).
$target = $null
for ($i = 2; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.Trim() -eq "") {
        $prev = $d.Paragraphs.Item($i - 1)
        if ($prev.Range.Text -like "Add new obstacles*") {
            $target = $cand
            break
        }
    }
}

$ratXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:r><w:t xml:space="preserve">Make Bat predict rat </w:t></w:r>' +
          '<w:proofErr w:type="gramStart"/>' +
          '<w:r><w:t>position</w:t></w:r>' +
          '<w:proofErr w:type="gramEnd"/>' +
          '</w:p>'
[void]$target.Range.InsertXML($ratXml)

# 2) Insert a new paragraph "Visual effect for cheese worth 1000 or more"
#    right after the "... (15%)" paragraph and before "Cheese size scales
#    with score".
$scoreIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*(15%)*") {
        $scoreIndex = $i
        break
    }
}

$scoreP = $d.Paragraphs.Item($scoreIndex)
$scoreP.Range.InsertParagraphAfter()
$newP = $d.Paragraphs.Item($scoreIndex + 1)
$newP.Range.Text = "Visual effect for cheese worth 1000 or more"
